$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new text would otherwise be auto-parsed by Excel as a pure number
# are pre-formatted as Text so the stored value keeps its exact original string form,
# matching the source data which stores prices/volumes as text.
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D50").NumberFormat = "@"

$ws.Range("D2").Value = "25.955.24"

$ws.Range("D3").Value = "1.637.12"
$ws.Range("E3").Value = "  -2.04%  "

$ws.Range("D4").Value = "1.003"
$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "209.19"
$ws.Range("E5").Value = "  -0.66%  "

$ws.Range("D6").Value = "0.5158"
$ws.Range("E6").Value = "  -1.80%  "

$ws.Range("E7").Value = "  +0.06%  "

$ws.Range("D8").Value = "0.2556"
$ws.Range("E8").Value = "  -3.68%  "

$ws.Range("D9").Value = "0.06225"
$ws.Range("E9").Value = "  -0.96%  "

$ws.Range("D10").Value = "20.26"
$ws.Range("E10").Value = "  -4.54%  "

$ws.Range("D11").Value = "0.07539"
$ws.Range("E11").Value = "  +0.14%  "

$ws.Range("D12").Value = "1.642.81"
$ws.Range("E12").Value = "  -1.27%  "

$ws.Range("D13").Value = "4.350"
$ws.Range("E13").Value = "  -2.16%  "

$ws.Range("D14").Value = "1.866.69"
$ws.Range("E14").Value = "  -1.71%  "

$ws.Range("D15").Value = "0.5387"
$ws.Range("E15").Value = "  -4.43%  "

$ws.Range("D16").Value = "0.0₅7909"
$ws.Range("E16").Value = "  -1.43%  "

$ws.Range("D17").Value = "64.70"
$ws.Range("E17").Value = "  -2.63%  "

$ws.Range("D18").Value = "25.975.23"

$ws.Range("E19").Value = "  -0.01%  "

$ws.Range("D20").Value = "4.628"
$ws.Range("E20").Value = "  -3.67%  "

$ws.Range("D21").Value = "184.73"
$ws.Range("E21").Value = "  -1.70%  "

$ws.Range("D22").Value = "9.966"
$ws.Range("E22").Value = "  -4.03%  "

$ws.Range("D23").Value = "6.071"
$ws.Range("E23").Value = "  -1.72%  "

$ws.Range("E24").Value = "  +0.01%  "

$ws.Range("D25").Value = "145.29"
$ws.Range("E25").Value = "  -2.03%  "

$ws.Range("D26").Value = "7.307"
$ws.Range("E26").Value = "  -3.81%  "

$ws.Range("D27").Value = "0.1187"
$ws.Range("E27").Value = "  -5.07%  "

$ws.Range("D28").Value = "15.45"
$ws.Range("E28").Value = "  -3.16%  "

$ws.Range("D29").Value = "1.373"
$ws.Range("E29").Value = "  +1.17%  "

$ws.Range("D30").Value = "0.05937"
$ws.Range("E30").Value = "  -4.82%  "

$ws.Range("E31").Value = "  -2.86%  "

$ws.Range("D32").Value = "3.338"
$ws.Range("E32").Value = "  -2.84%  "

$ws.Range("D33").Value = "3.346"
$ws.Range("E33").Value = "  -3.86%  "

$ws.Range("D34").Value = "1.600"
$ws.Range("E34").Value = "  -1.80%  "

$ws.Range("D35").Value = "0.9665"
$ws.Range("E35").Value = "  -3.25%  "

$ws.Range("D36").Value = "2.382"
$ws.Range("E36").Value = "  -0.94%  "

$ws.Range("D37").Value = "2.733"
$ws.Range("E37").Value = "  +0.66%  "

$ws.Range("D38").Value = "0.5829"
$ws.Range("E38").Value = "  -3.57%  "

$ws.Range("D39").Value = "0.01592"
$ws.Range("E39").Value = "  -1.32%  "

$ws.Range("D40").Value = "1.002"
$ws.Range("E40").Value = "  -0.25%  "

$ws.Range("B41").Value = "TrustWalletToken"
$ws.Range("C41").Value = "https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt"
$ws.Range("D41").Value = "0.8403"
$ws.Range("E41").Value = "  -3.02%  "

$ws.Range("B42").Value = "Maker"
$ws.Range("C42").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D42").Value = "1.038.57"
$ws.Range("E42").Value = "  -3.53%  "

$ws.Range("D43").Value = "5.681"
$ws.Range("E43").Value = "  -7.10%  "

$ws.Range("D44").Value = "99.57"
$ws.Range("E44").Value = "  -0.45%  "

$ws.Range("D45").Value = "1.792.54"
$ws.Range("E45").Value = "  -1.58%  "

$ws.Range("D46").Value = "0.0₈106"
$ws.Range("E46").Value = "  -2.47%  "

$ws.Range("D47").Value = "1.001"
$ws.Range("E47").Value = "  +0.05%  "

$ws.Range("D48").Value = "54.24"
$ws.Range("E48").Value = "  -3.40%  "

$ws.Range("D49").Value = "7.916"
$ws.Range("E49").Value = "  -1.18%  "

$ws.Range("D50").Value = "0.05188"
$ws.Range("E50").Value = "  -1.01%  "

$ws.Range("E51").Value = "  -0.54%  "
